$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.761.81"
$ws.Range("E2").Value = "  +7.14%  "

$ws.Range("D3").Value = "1.810.24"
$ws.Range("E3").Value = "  +4.53%  "

$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "'252.56"
$ws.Range("E5").Value = "  +4.25%  "

$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").Value = "'0.4958"
$ws.Range("E7").Value = "  +0.92%  "

$ws.Range("D8").Value = "'0.2805"
$ws.Range("E8").Value = "  +7.50%  "

$ws.Range("D9").Value = "'0.06394"
$ws.Range("E9").Value = "  +2.75%  "

$ws.Range("D10").Value = "1.807.22"
$ws.Range("E10").Value = "  +3.98%  "

$ws.Range("D11").Value = "'16.84"
$ws.Range("E11").Value = "  +4.70%  "

$ws.Range("D12").Value = "'0.07109"
$ws.Range("E12").Value = "  +2.84%  "

$ws.Range("D13").Value = "'0.6472"
$ws.Range("E13").Value = "  +5.74%  "

$ws.Range("D14").Value = "'4.703"
$ws.Range("E14").Value = "  +4.22%  "

$ws.Range("D15").Value = "'82.26"
$ws.Range("E15").Value = "  +6.33%  "

$ws.Range("D16").Value = "28.747.06"
$ws.Range("E16").Value = "  +7.98%  "

$ws.Range("D17").Value = "'0.9996"
$ws.Range("E17").Value = "  +0.10%  "

$ws.Range("D18").Value = "'0.000007365"
$ws.Range("E18").Value = "  +2.38%  "

$ws.Range("D19").Value = "'0.9992"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").Value = "'12.30"
$ws.Range("E20").Value = "  +7.08%  "

$ws.Range("D21").Value = "2.041.20"
$ws.Range("E21").Value = "  +4.17%  "

$ws.Range("D22").Value = "'4.612"
$ws.Range("E22").Value = "  +3.83%  "

$ws.Range("D23").Value = "'8.882"
$ws.Range("E23").Value = "  +3.54%  "

$ws.Range("D24").Value = "'5.303"
$ws.Range("E24").Value = "  +3.37%  "

$ws.Range("D25").Value = "'142.93"
$ws.Range("E25").Value = "  +2.93%  "

$ws.Range("D26").Value = "'16.04"
$ws.Range("E26").Value = "  +4.70%  "

$ws.Range("D27").Value = "'1.883"
$ws.Range("E27").Value = "  +5.29%  "

$ws.Range("E28").Value = "  +5.16%  "

$ws.Range("D29").Value = "'1.387"
$ws.Range("E29").Value = "  +0.53%  "

$ws.Range("D30").Value = "'4.185"
$ws.Range("E30").Value = "  +6.13%  "

$ws.Range("D31").Value = "'0.08363"
$ws.Range("E31").Value = "  +4.53%  "

$ws.Range("D32").Value = "'3.853"
$ws.Range("E32").Value = "  +4.52%  "

$ws.Range("D33").Value = "'0.04951"
$ws.Range("E33").Value = "  +9.25%  "

$ws.Range("E34").Value = "  +8.51%  "

$ws.Range("D35").Value = "'0.6715"
$ws.Range("E35").Value = "  +7.44%  "

$ws.Range("D36").Value = "'2.665"
$ws.Range("E36").Value = "  +2.36%  "

$ws.Range("D37").Value = "'0.9641"
$ws.Range("E37").Value = "  +2.53%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.636"
$ws.Range("E38").Value = "  +7.68%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'2.151"
$ws.Range("E39").Value = "  +4.84%  "

$ws.Range("D40").Value = "'0.01604"
$ws.Range("E40").Value = "  +6.51%  "

$ws.Range("D41").Value = "'5.965"
$ws.Range("E41").Value = "  +5.49%  "

$ws.Range("D42").Value = "'0.9996"
$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'101.36"
$ws.Range("E43").Value = "  +1.74%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.4128"
$ws.Range("E44").Value = "  +6.70%  "

$ws.Range("D45").Value = "'7.235"
$ws.Range("E45").Value = "  +4.11%  "

$ws.Range("D46").Value = "'0.1225"
$ws.Range("E46").Value = "  +5.37%  "

$ws.Range("D47").Value = "'0.05492"
$ws.Range("E47").Value = "  +1.90%  "

$ws.Range("D48").Value = "'8.207"
$ws.Range("E48").Value = "  +3.02%  "

$ws.Range("D49").Value = "'31.32"
$ws.Range("E49").Value = "  +3.53%  "

$ws.Range("D50").Value = "'0.3620"
$ws.Range("E50").Value = "  +6.69%  "

$ws.Range("D51").Value = "'1.304"
$ws.Range("E51").Value = "  +4.69%  "
